# Auto-generated edit script applying cached-value updates to Garuda_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 40: H40,I40,J40,K40,L40,M40,N40
$ws.Range("H40").Value = 1881.375
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1827.4546
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1827.4546
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2177.4546
# row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 6382.4116
$ws.Range("I74").Value = 12800.2
$ws.Range("J74").Value = 3708.3333
$ws.Range("K74").Value = 12800.2
$ws.Range("L74").Value = 3708.3333
$ws.Range("M74").Value = -11864.2
$ws.Range("N74").Value = -5580.3333
# row 76: H76,I76,J76,K76,L76,M76,N76
$ws.Range("H76").Value = 69533.53
$ws.Range("I76").Value = 79538.69500000001
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 79538.69500000001
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -79223.69500000001
$ws.Range("N76").Value = -5130
# row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 6382.4116
$ws.Range("I77").Value = 12800.2
$ws.Range("J77").Value = 3708.3333
$ws.Range("K77").Value = 64001
$ws.Range("L77").Value = 18541.6665
$ws.Range("M77").Value = -59321
$ws.Range("N77").Value = -27901.6665
# row 79: H79,I79,J79,K79,L79,M79,N79
$ws.Range("H79").Value = 69533.53
$ws.Range("I79").Value = 79538.69500000001
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 79538.69500000001
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -78446.69500000001
$ws.Range("N79").Value = -6684
# row 86: H86,J86,L86,N86
$ws.Range("H86").Value = 54155.21
$ws.Range("J86").Value = 2105
$ws.Range("L86").Value = 2105
$ws.Range("N86").Value = -4351
# row 89: H89,J89,L89,N89
$ws.Range("H89").Value = 54155.21
$ws.Range("J89").Value = 2105
$ws.Range("L89").Value = 10525
$ws.Range("N89").Value = -21757
# row 98: H98,I98,J98,K98,L98,M98,N98
$ws.Range("H98").Value = 38819.344
$ws.Range("I98").Value = 5035.625
$ws.Range("J98").Value = 200981.2
$ws.Range("K98").Value = 5035.625
$ws.Range("L98").Value = 200981.2
$ws.Range("M98").Value = -3537.625
$ws.Range("N98").Value = -203977.2
# row 106: H106,I106,J106,K106,L106,M106,N106
$ws.Range("H106").Value = 50100580
$ws.Range("I106").Value = 111756.664
$ws.Range("J106").Value = 500000000
$ws.Range("K106").Value = 111756.664
$ws.Range("L106").Value = 500000000
$ws.Range("M106").Value = -111125.664
$ws.Range("N106").Value = -500001262
# row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 38819.344
$ws.Range("I122").Value = 5035.625
$ws.Range("J122").Value = 200981.2
$ws.Range("K122").Value = 15106.875
$ws.Range("L122").Value = 602943.6000000001
$ws.Range("M122").Value = -12656.875
$ws.Range("N122").Value = -607843.6000000001
# row 135: H135,I135,J135,K135,L135,M135,N135
$ws.Range("H135").Value = 566.0851
$ws.Range("I135").Value = 513.1739
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 4618.5651
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -2083.5651
$ws.Range("N135").Value = -32070
# row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 1365.8776
$ws.Range("I137").Value = 1121.7368
$ws.Range("J137").Value = 2209.2727
$ws.Range("K137").Value = 3365.2104
$ws.Range("L137").Value = 6627.8181
$ws.Range("M137").Value = -815.2103999999999
$ws.Range("N137").Value = -11727.8181

$ws = $wb.Worksheets.Item("ARM")
# row 32: H32,I32,J32,K32,L32,M32,N32
$ws.Range("H32").Value = 14119.883
$ws.Range("I32").Value = 15188.94
$ws.Range("J32").Value = 5139.8
$ws.Range("K32").Value = 15188.94
$ws.Range("L32").Value = 5139.8
$ws.Range("M32").Value = -14901.94
$ws.Range("N32").Value = -5713.8
# row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 1507.0256
$ws.Range("I61").Value = 958.9231
$ws.Range("J61").Value = 2603.2307
$ws.Range("K61").Value = 958.9231
$ws.Range("L61").Value = 2603.2307
$ws.Range("M61").Value = -746.9231
$ws.Range("N61").Value = -3027.2307
# row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 1181.5333
$ws.Range("I74").Value = 1080.875
$ws.Range("J74").Value = 1584.1666
$ws.Range("K74").Value = 1080.875
$ws.Range("L74").Value = 1584.1666
$ws.Range("M74").Value = -206.875
$ws.Range("N74").Value = -3332.1666
# row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 1181.5333
$ws.Range("I77").Value = 1080.875
$ws.Range("J77").Value = 1584.1666
$ws.Range("K77").Value = 5404.375
$ws.Range("L77").Value = 7920.833000000001
$ws.Range("M77").Value = -1036.375
$ws.Range("N77").Value = -16656.833
# row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 3632.0876
$ws.Range("I132").Value = 4671.242
$ws.Range("J132").Value = 2203.25
$ws.Range("K132").Value = 14013.726
$ws.Range("L132").Value = 6609.75
$ws.Range("M132").Value = -11483.726
$ws.Range("N132").Value = -11669.75
# row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 1507.0256
$ws.Range("I136").Value = 958.9231
$ws.Range("J136").Value = 2603.2307
$ws.Range("K136").Value = 2876.7693
$ws.Range("L136").Value = 7809.6921
$ws.Range("M136").Value = -326.7692999999999
$ws.Range("N136").Value = -12909.6921

$ws = $wb.Worksheets.Item("BSM")
# row 105: H105,I105,K105,M105
$ws.Range("H105").Value = 2835.1052
$ws.Range("I105").Value = 2402.647
$ws.Range("K105").Value = 2402.647
$ws.Range("M105").Value = -655.6469999999999
# row 113: H113,I113,K113,M113
$ws.Range("H113").Value = 28385
$ws.Range("I113").Value = 28385
$ws.Range("K113").Value = 28385
$ws.Range("M113").Value = -26215

$ws = $wb.Worksheets.Item("CRP")
# row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 5955168
$ws.Range("I31").Value = 1918.5
$ws.Range("J31").Value = 20838292
$ws.Range("K31").Value = 1918.5
$ws.Range("L31").Value = 20838292
$ws.Range("M31").Value = -1623.5
$ws.Range("N31").Value = -20838882
# row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 5955168
$ws.Range("I34").Value = 1918.5
$ws.Range("J34").Value = 20838292
$ws.Range("K34").Value = 1918.5
$ws.Range("L34").Value = 20838292
$ws.Range("M34").Value = -1716.5
$ws.Range("N34").Value = -20838696
# row 58: H58,I58,J58,K58,L58,M58,N58
$ws.Range("H58").Value = 852.43335
$ws.Range("I58").Value = 993.3333
$ws.Range("J58").Value = 641.0833
$ws.Range("K58").Value = 993.3333
$ws.Range("L58").Value = 641.0833
$ws.Range("M58").Value = -790.3333
$ws.Range("N58").Value = -1047.0833
# row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1661.7587
$ws.Range("I132").Value = 1647.075
$ws.Range("J132").Value = 1694.3889
$ws.Range("K132").Value = 4941.225
$ws.Range("L132").Value = 5083.1667
$ws.Range("M132").Value = -2411.225
$ws.Range("N132").Value = -10143.1667
# row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 1165.4849
$ws.Range("I134").Value = 1115.7
$ws.Range("J134").Value = 1663.3334
$ws.Range("K134").Value = 3347.1
$ws.Range("L134").Value = 4990.0002
$ws.Range("M134").Value = -812.1000000000004
$ws.Range("N134").Value = -10060.0002
# row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 852.43335
$ws.Range("I136").Value = 993.3333
$ws.Range("J136").Value = 641.0833
$ws.Range("K136").Value = 2979.9999
$ws.Range("L136").Value = 1923.2499
$ws.Range("M136").Value = -429.9998999999998
$ws.Range("N136").Value = -7023.2499

$ws = $wb.Worksheets.Item("CUL")
# row 131: H131,I131,J131,K131,L131,M131,N131
$ws.Range("H131").Value = 5535.4546
$ws.Range("I131").Value = 6587.778
$ws.Range("J131").Value = 800
$ws.Range("K131").Value = 19763.334
$ws.Range("L131").Value = 2400
$ws.Range("M131").Value = -14723.334
$ws.Range("N131").Value = -12480

$ws = $wb.Worksheets.Item("GSM")
# row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Range("H70").Value = 13425171
$ws.Range("I70").Value = 15004130
$ws.Range("J70").Value = 4025
$ws.Range("K70").Value = 15004130
$ws.Range("L70").Value = 4025
$ws.Range("M70").Value = -15003860
$ws.Range("N70").Value = -4565
# row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Range("H73").Value = 13425171
$ws.Range("I73").Value = 15004130
$ws.Range("J73").Value = 4025
$ws.Range("K73").Value = 15004130
$ws.Range("L73").Value = 4025
$ws.Range("M73").Value = -15003194
$ws.Range("N73").Value = -5897
# row 97: H97,J97,L97,N97
$ws.Range("H97").Value = 748.6842
$ws.Range("J97").Value = 940
$ws.Range("L97").Value = 940
$ws.Range("N97").Value = -1932
# row 131: H131,J131,L131,N131
$ws.Range("H131").Value = 27334
$ws.Range("J131").Value = 27334
$ws.Range("L131").Value = 27334
$ws.Range("N131").Value = -37414
# row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 19978.11
$ws.Range("I132").Value = 26680.2
$ws.Range("J132").Value = 2105.8667
$ws.Range("K132").Value = 80040.60000000001
$ws.Range("L132").Value = 6317.6001
$ws.Range("M132").Value = -77510.60000000001
$ws.Range("N132").Value = -11377.6001

$ws = $wb.Worksheets.Item("LTW")
# row 7: H7,I7,J7,K7,L7,M7,N7
$ws.Range("H7").Value = 16668190
$ws.Range("I7").Value = 1195.3158
$ws.Range("J7").Value = 45456636
$ws.Range("K7").Value = 1195.3158
$ws.Range("L7").Value = 45456636
$ws.Range("M7").Value = -1083.3158
$ws.Range("N7").Value = -45456860
# row 21: H21,J21,L21,N21
$ws.Range("H21").Value = 12000
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12348
# row 62: H62,J62,L62,N62
$ws.Range("H62").Value = 29000
$ws.Range("J62").Value = 29000
$ws.Range("L62").Value = 29000
$ws.Range("N62").Value = -30248
# row 65: H65,J65,L65,N65
$ws.Range("H65").Value = 29000
$ws.Range("J65").Value = 29000
$ws.Range("L65").Value = 87000
$ws.Range("N65").Value = -93240
# row 68: H68,I68,J68,K68,L68,M68,N68
$ws.Range("H68").Value = 1550.1428
$ws.Range("I68").Value = 1137.75
$ws.Range("J68").Value = 2100
$ws.Range("K68").Value = 1137.75
$ws.Range("L68").Value = 2100
$ws.Range("M68").Value = -388.75
$ws.Range("N68").Value = -3598
# row 71: H71,I71,J71,K71,L71,M71,N71
$ws.Range("H71").Value = 1550.1428
$ws.Range("I71").Value = 1137.75
$ws.Range("J71").Value = 2100
$ws.Range("K71").Value = 5688.75
$ws.Range("L71").Value = 10500
$ws.Range("M71").Value = -1944.75
$ws.Range("N71").Value = -17988
# row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 16668190
$ws.Range("I126").Value = 1195.3158
$ws.Range("J126").Value = 45456636
$ws.Range("K126").Value = 3585.9474
$ws.Range("L126").Value = 136369908
$ws.Range("M126").Value = -1115.9474
$ws.Range("N126").Value = -136374848
# row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 5986.9556
$ws.Range("I132").Value = 8533.179
$ws.Range("J132").Value = 1793.1765
$ws.Range("K132").Value = 25599.537
$ws.Range("L132").Value = 5379.529500000001
$ws.Range("M132").Value = -23069.537
$ws.Range("N132").Value = -10439.5295
# row 133: H133,J133,L133,N133
$ws.Range("H133").Value = 19064.125
$ws.Range("J133").Value = 19064.125
$ws.Range("L133").Value = 19064.125
$ws.Range("N133").Value = -24124.125
# row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 3755.7437
$ws.Range("I136").Value = 4171.517
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 12514.551
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -9964.550999999999
$ws.Range("N136").Value = -12750

$ws = $wb.Worksheets.Item("WVR")
# row 8: H8,J8,L8,N8
$ws.Range("H8").Value = 1600
$ws.Range("J8").Value = 1600
$ws.Range("L8").Value = 1600
$ws.Range("N8").Value = -1880
# row 82: H82,I82,K82,M82
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2617
# row 85: H85,I85,K85,M85
$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1674
# row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 858.0492
$ws.Range("I132").Value = 789.7778
$ws.Range("K132").Value = 2369.3334
$ws.Range("M132").Value = 160.6666
# row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 3104.9075
$ws.Range("I136").Value = 3317.8262
$ws.Range("K136").Value = 9953.4786
$ws.Range("M136").Value = -7403.4786
